$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price column D, Volume(1h) column E)
# D-column numeric-looking values are written with NumberFormat "@" (Text)
# so the engine stores them as text, matching the source inline strings
# instead of auto-converting them to numbers.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '63.142.85'

$ws.Range('E2').Value = '  -0.04%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.574.53'

$ws.Range('E3').Value = '  +0.23%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '591.74'

$ws.Range('E5').Value = '  +1.21%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '144.47'

$ws.Range('E6').Value = '  -2.22%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.590'

$ws.Range('E8').Value = '  -1.82%  '

$ws.Range('E9').Value = '  -2.10%  '

$ws.Range('E10').Value = '  -1.01%  '

$ws.Range('E11').Value = '  -0.20%  '

$ws.Range('E12').Value = '  -1.65%  '

$ws.Range('E13').Value = '  -0.89%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.039.27'

$ws.Range('E14').Value = '  +0.34%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '63.050.82'

$ws.Range('E15').Value = '  -0.11%  '

$ws.Range('E16').Value = '  -0.99%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.580.68'

$ws.Range('E17').Value = '  +0.32%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '11.07'

$ws.Range('E18').Value = '  -2.50%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '340.99'

$ws.Range('E19').Value = '  -0.83%  '

$ws.Range('E20').Value = '  -1.90%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.64'

$ws.Range('E21').Value = '  -3.69%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.75'

$ws.Range('E23').Value = '  +3.76%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '67.78'

$ws.Range('E24').Value = '  +1.33%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.58'

$ws.Range('E25').Value = '  +6.91%  '

$ws.Range('E26').Value = '  -0.85%  '

$ws.Range('E27').Value = '  -3.14%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.999'

$ws.Range('E28').Value = '  -0.01%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.91'

$ws.Range('E29').Value = '  -3.08%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.25'

$ws.Range('E30').Value = '  -2.70%  '

$ws.Range('E31').Value = '  -1.87%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '468.56'

$ws.Range('E32').Value = '  +0.57%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0₃0801'

$ws.Range('E34').Value = '  +3.14%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '176.52'

$ws.Range('E35').Value = '  +0.48%  '

$ws.Range('E36').Value = '  +0.04%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.400'

$ws.Range('E37').Value = '  -1.91%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '18.85'

$ws.Range('E38').Value = '  -1.94%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '4.56'

$ws.Range('E39').Value = '  -0.13%  '

$ws.Range('E40').Value = '  +0.00%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.69'

$ws.Range('E41').Value = '  -3.42%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '40.05'

$ws.Range('E42').Value = '  +1.06%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '158.00'

$ws.Range('E43').Value = '  +4.11%  '

$ws.Range('E44').Value = '  -3.29%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '21.36'

$ws.Range('E45').Value = '  +1.65%  '

$ws.Range('E46').Value = '  +3.40%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0538'

$ws.Range('E47').Value = '  -1.65%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0962'

$ws.Range('E48').Value = '  -1.61%  '

$ws.Range('E49').Value = '  -1.23%  '

$ws.Range('E50').Value = '  -1.91%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '11.39'

$ws.Range('E51').Value = '  +0.05%  '
